# Update the two "Output: Points" formulas in J11 and J12 so that the
# points-per-meter multiplier reflects the Large Hill constant (1.8)
# instead of the previously hard-coded 2, and have J12 reference the
# "Output: Distance" cell (I12) like its sibling rows do.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J11").Formula = "=60+(H11-F11)*1.8"
$ws.Range("J12").Formula = "=60+(I12-120)*1.8"

# Leave the selection on J12, matching where the edits were last made.
$null = $ws.Range("J12").Select()
